$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 124 currently holds the last data row (2020-10-05 label, values for that date).
# Insert a new row before it containing the data formerly attributed to "2020-10-05"
# re-labelled data, and push the old row 124 values down the sequence as Oct 6-7 data,
# then append Oct 1 - Oct 7 daily rows plus a final "2020-10-08" summary row.

# Update row 124 (previously the final row) to the new Oct-01 raw values (text label stays as-is, shared string updated elsewhere)
$ws.Range("B124").Value = 748315
$ws.Range("C124").Value = 884896
$ws.Range("D124").Value = 88659
$ws.Range("E124").Value = 78078
$ws.Range("F124").Value = 23.88

# New daily rows 125-131, each with a real date value formatted as a date
$dates = @(44105,44106,44107,44108,44109,44110,44111)
$bvals = @(748315,753090,757953,761665,765082,769558,774020)
$cvals = @(884896,893324,901110,907331,913155,933316,940994)
$dvals = @(88659,89183,90194,85743,80345,32797,36802)
$evals = @(78078,78492,78880,79088,79268,79714,80083)
$fvals = @(23.88031778061378,23.810832702598631,23.7492298335121,23.707535465066659,23.675501449517832,23.652018431359298,23.618123562698639)

for ($i = 0; $i -lt 7; $i++) {
    $r = 125 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 1).NumberFormat = "mm-dd-yy"
    $ws.Cells.Item($r, 2).Value = $bvals[$i]
    $ws.Cells.Item($r, 3).Value = $cvals[$i]
    $ws.Cells.Item($r, 4).Value = $dvals[$i]
    $ws.Cells.Item($r, 5).Value = $evals[$i]
    $ws.Cells.Item($r, 6).Value = $fvals[$i]
}

# Final summary row 132 labelled "2020-10-08"
$ws.Range("A132").Value = "2020-10-08"
$ws.Range("B132").Value = 804488
$ws.Range("C132").Value = 948928
$ws.Range("D132").Value = 299866
$ws.Range("E132").Value = 83096
$ws.Range("F132").Value = 23.5
